$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(104).Insert()

$ws.Cells.Item(104, 1).Value = 10
$ws.Cells.Item(104, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(104, 3).Value = "La Araucanía"
$ws.Cells.Item(104, 4).Value = 45033
$ws.Cells.Item(104, 5).Value = 9
$ws.Cells.Item(104, 6).Value = "Fruta"
$ws.Cells.Item(104, 7).Value = 100104
$ws.Cells.Item(104, 8).Value = "Frutos de pepita"
$ws.Cells.Item(104, 9).Value = 100104001
$ws.Cells.Item(104, 10).Value = "Granada"
$ws.Cells.Item(104, 11).Value = "Wonderfull"
$ws.Cells.Item(104, 12).Value = "Primera"
$ws.Cells.Item(104, 13).Value = 150
$ws.Cells.Item(104, 14).Value = 20000
$ws.Cells.Item(104, 15).Value = 20000
$ws.Cells.Item(104, 16).Value = 20000
$ws.Cells.Item(104, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(104, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(104, 19).Value = 1333
$ws.Cells.Item(104, 20).Value = 15
